# "pic defined on admission day"
# Shift the weekly admission dates by a few days and refresh the
# dependent parameter table (mpic/H/I columns etc.) on the `params` sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Admission dates (column A, rows 3-5) -------------------------------
# 07.03.2020 -> 13.03.2020
# 10.03.2020 -> 18.03.2020
# 21.03.2020 -> 23.03.2020
$ws.Range("A3").Value = "13.03.2020"
$ws.Range("A4").Value = "18.03.2020"
$ws.Range("A5").Value = "23.03.2020"

# --- Row 2 ---------------------------------------------------------------
$ws.Range("C2").Value = 0.05
$ws.Range("H2").Value = 10
$ws.Range("I2").Value = 20

# --- Row 3 ---------------------------------------------------------------
$ws.Range("C3").Value = 0.05
$ws.Range("D3").Value = 0.3
$ws.Range("H3").Value = 10
$ws.Range("I3").Value = 20

# --- Row 4 ---------------------------------------------------------------
$ws.Range("C4").Value = 0.05
$ws.Range("D4").Value = 0.01
$ws.Range("H4").Value = 10
$ws.Range("I4").Value = 20

# --- Row 5 ---------------------------------------------------------------
$ws.Range("C5").Value = 0.05
$ws.Range("H5").Value = 10
$ws.Range("I5").Value = 20

# --- Selected cell moved from A6 to F6 -----------------------------------
$null = $ws.Range("F6").Select()
